# Add a new "3segment" worksheet at the end of the workbook (after the
# current last sheet, "MPrior") and populate it with the theta1/theta2/MP
# lookup table, matching the author's "basic working version with angle
# pdf" commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "3segment"

# Header row
$ws.Range("B1").Value = "theta 1"
$ws.Range("C1").Value = "theta 2"
$ws.Range("D1").Value = "MP"

# Data rows
$ws.Range("A2").Value = "D"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "(0,1)"

$ws.Range("A3").Value = "P"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "(0.5,1)"

$ws.Range("A4").Value = "T"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "(0.5,0)"

$ws.Range("A5").Value = "X"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "(0.5,0.5)"

$ws.Range("A6").Value = "Y"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "(1,0.5)"

$ws.Range("A7").Value = "Z"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "(nan,nan)"

# The newly added sheet becomes the selected/active tab, as in the diff.
[void]$ws.Select()
[void]$ws.Range("D15").Select()
